# Update TPM-derived values in the LR-pairs sheet (Gpi1-Amfr)
# per "update scripts wuth new tpm" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 19.71326633333333
$ws.Range("H2").Value = 59.139799
$ws.Range("I2").Value = 0.2311148647321797
$ws.Range("J2").Value = 0.2311148647321797
$ws.Range("M2").Value = 10.44283133333333
$ws.Range("N2").Value = 31.328494
$ws.Range("O2").Value = 0.2278338596647878
$ws.Range("P2").Value = 0.2278338596647878
$ws.Range("Q2").Value = 205.8623153480784
$ws.Range("R2").Value = 1852.760838132706
$ws.Range("S2").Value = 0.05265579165783784
$ws.Range("T2").Value = 0.05265579165783785

# Row 3
$ws.Range("G3").Value = 19.71326633333333
$ws.Range("H3").Value = 59.139799
$ws.Range("I3").Value = 0.2311148647321797
$ws.Range("J3").Value = 0.2311148647321797
$ws.Range("O3").Value = 0.5060228699757219
$ws.Range("P3").Value = 0.5060228699757219
$ws.Range("Q3").Value = 457.2236970639426
$ws.Range("R3").Value = 4115.013273575483
$ws.Range("S3").Value = 0.1169494071458283
$ws.Range("T3").Value = 0.1169494071458283

# Row 4
$ws.Range("G4").Value = 19.71326633333333
$ws.Range("H4").Value = 59.139799
$ws.Range("I4").Value = 0.2311148647321797
$ws.Range("J4").Value = 0.2311148647321797
$ws.Range("O4").Value = 0.2661432703594902
$ws.Range("P4").Value = 0.2661432703594903
$ws.Range("Q4").Value = 240.4772931079041
$ws.Range("R4").Value = 2164.295637971137
$ws.Range("S4").Value = 0.06150966592851351
$ws.Range("T4").Value = 0.06150966592851351

# Row 5
$ws.Range("I5").Value = 0.3593964817703036
$ws.Range("J5").Value = 0.3593964817703036
$ws.Range("M5").Value = 10.44283133333333
$ws.Range("N5").Value = 31.328494
$ws.Range("O5").Value = 0.2278338596647878
$ws.Range("P5").Value = 0.2278338596647878
$ws.Range("Q5").Value = 320.1273615651021
$ws.Range("R5").Value = 2881.14625408592
$ws.Range("S5").Value = 0.08188268759167383
$ws.Range("T5").Value = 0.08188268759167384

# Row 6
$ws.Range("I6").Value = 0.3593964817703036
$ws.Range("J6").Value = 0.3593964817703036
$ws.Range("O6").Value = 0.5060228699757219
$ws.Range("P6").Value = 0.5060228699757219
$ws.Range("S6").Value = 0.1818628391645863
$ws.Range("T6").Value = 0.1818628391645863

# Row 7
$ws.Range("I7").Value = 0.3593964817703036
$ws.Range("J7").Value = 0.3593964817703036
$ws.Range("O7").Value = 0.2661432703594902
$ws.Range("P7").Value = 0.2661432703594903
$ws.Range("S7").Value = 0.09565095501404351
$ws.Range("T7").Value = 0.09565095501404354

# Row 8
$ws.Range("I8").Value = 0.4094886534975166
$ws.Range("J8").Value = 0.4094886534975166
$ws.Range("M8").Value = 10.44283133333333
$ws.Range("N8").Value = 31.328494
$ws.Range("O8").Value = 0.2278338596647878
$ws.Range("P8").Value = 0.2278338596647878
$ws.Range("Q8").Value = 364.7462590320715
$ws.Range("R8").Value = 3282.716331288644
$ws.Range("S8").Value = 0.09329538041527612
$ws.Range("T8").Value = 0.09329538041527613

# Row 9
$ws.Range("I9").Value = 0.4094886534975166
$ws.Range("J9").Value = 0.4094886534975166
$ws.Range("O9").Value = 0.5060228699757219
$ws.Range("P9").Value = 0.5060228699757219
$ws.Range("S9").Value = 0.2072106236653073
$ws.Range("T9").Value = 0.2072106236653073

# Row 10
$ws.Range("I10").Value = 0.4094886534975166
$ws.Range("J10").Value = 0.4094886534975166
$ws.Range("O10").Value = 0.2661432703594902
$ws.Range("P10").Value = 0.2661432703594903
$ws.Range("S10").Value = 0.1089826494169332
$ws.Range("T10").Value = 0.1089826494169332
